$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting existing columns B:H to C:I
$ws.Range("B1").EntireColumn.Insert()

# Copy the header formatting from the (now shifted) "id" header cell onto the
# new column's header cell, then set its text.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B1").Value = "env"
